# Update Pooh Points site
# Applies the 2026-01-27 SEC ByOwner refresh:
#  - Players sheet: updated live stat snapshot for rows 2-33 (new game
#    clock/status plus refreshed pooh/pts/reb/ast/stl/blk/to/min numbers),
#    a few rows where two players on the same owner/team swapped places,
#    and one brand-new waiver-wire row (34) for Sebastian Mack.
#  - OwnerTotals sheet: re-sorted/recomputed starter totals that follow
#    from the Players sheet refresh.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Columns on the Players sheet: A date, B owner, C started_today, D player,
# E team, F game, G status, H pooh, I pts, J reb, K ast, L stl, M blk, N to, O min.
# Column A (date) is identical before/after for every row, so we leave it
# untouched here (writing literal "2026-01-27" strings through .Value would
# otherwise get auto-coerced into an Excel date serial number).
$rows = @(
  @(2,  "Booz",      "Yes", "Labaron Philon Jr.",   "ALA", "MIZ@ALA", "9:52 - 1st Half",  13, 11, 1,  5, 1, 0, 1, 17),
  @(3,  "Booz",      "No",  "Derrion Reid",          "OU",  "ARK@OU",  "11:16 - 2nd Half", 11, 12, 5,  0, 0, 0, 2, 29),
  @(4,  "Booz",      "No",  "Billy Richmond III",    "ARK", "ARK@OU",  "11:16 - 2nd Half", 9,  9,  2,  0, 2, 0, 0, 15),
  @(5,  "Booz",      "No",  "Jacob Crews",           "MIZ", "MIZ@ALA", "9:52 - 1st Half",  0,  0,  3,  0, 0, 0, 0, 15),
  @(6,  "CDL",       "Yes", "Trevon Brazile",        "ARK", "ARK@OU",  "11:16 - 2nd Half", 13, 6,  7,  4, 0, 1, 2, 31),
  @(7,  "CDL",       "No",  "Nijel Pack",            "OU",  "ARK@OU",  "11:16 - 2nd Half", 16, 19, 1,  2, 0, 0, 1, 28),
  @(8,  "CDL",       "No",  "Meleek Thomas",         "ARK", "ARK@OU",  "11:16 - 2nd Half", 10, 14, 3,  1, 0, 0, 0, 30),
  @(9,  "Clay",      "Yes", "Mark Mitchell",         "MIZ", "MIZ@ALA", "9:52 - 1st Half",  5,  4,  4,  2, 1, 0, 0, 18),
  @(10, "Clay",      "Yes", "Aden Holloway",         "ALA", "MIZ@ALA", "9:52 - 1st Half",  4,  3,  1,  4, 0, 0, 1, 11),
  @(11, "Clay",      "No",  "Kirill Elatontsev",     "OU",  "ARK@OU",  "11:16 - 2nd Half", 6,  2,  3,  0, 1, 0, 0, 12),
  @(12, "Clay",      "No",  "Taylor Bol Bowen",      "ALA", "MIZ@ALA", "9:52 - 1st Half",  3,  2,  6,  0, 0, 0, 1, 11),
  @(13, "Hal",       "Yes", "Jayden Stone",          "MIZ", "MIZ@ALA", "9:52 - 1st Half",  6,  6,  3,  1, 1, 0, 1, 16),
  @(14, "Hal",       "No",  "Nick Pringle",          "ARK", "ARK@OU",  "11:16 - 2nd Half", 11, 8,  2,  1, 1, 0, 1, 21),
  @(15, "Hal",       "No",  "Jadon Jones",           "OU",  "ARK@OU",  "11:16 - 2nd Half", 4,  3,  0,  1, 0, 1, 0, 9),
  @(16, "Hal",       "No",  "Houston Mallette",      "ALA", "MIZ@ALA", "9:52 - 1st Half",  1,  0,  4,  0, 0, 0, 0, 13),
  @(17, "Hal",       "No",  "Anthony Robinson II",   "MIZ", "MIZ@ALA", "9:52 - 1st Half",  -1, 3,  0,  0, 0, 0, 1, 8),
  @(18, "Mark",      "Yes", "Darius Acuff Jr.",      "ARK", "ARK@OU",  "11:16 - 2nd Half", 21, 17, 2,  9, 0, 0, 0, 32),
  @(19, "Mark",      "No",  "Xzayvier Brown",        "OU",  "ARK@OU",  "11:16 - 2nd Half", 11, 11, 6,  2, 0, 0, 1, 34),
  @(20, "Mark",      "No",  "Jeff Nwankwo",          "OU",  "ARK@OU",  "11:16 - 2nd Half", 3,  4,  1,  0, 0, 0, 0, 13),
  @(21, "Ron",       "No",  "Mohamed Wague",         "OU",  "ARK@OU",  "11:16 - 2nd Half", 14, 7,  10, 4, 1, 1, 4, 24),
  @(22, "Ron",       "No",  "Latrell Wrightsell",    "ALA", "MIZ@ALA", "9:52 - 1st Half",  8,  6,  2,  2, 1, 0, 1, 17),
  @(23, "Tar",       "Yes", "Tae Davis",             "OU",  "ARK@OU",  "11:16 - 2nd Half", 11, 14, 3,  2, 0, 0, 2, 31),
  @(24, "Tar",       "Yes", "Aiden Sherrell",        "ALA", "MIZ@ALA", "9:52 - 1st Half",  7,  8,  1,  0, 0, 0, 0, 13),
  @(25, "Tar",       "No",  "Karter Knox",           "ARK", "ARK@OU",  "11:16 - 2nd Half", 12, 11, 2,  0, 1, 1, 0, 25),
  @(26, "Undrafted", "No",  "Malique Ewin",          "ARK", "ARK@OU",  "11:16 - 2nd Half", 12, 10, 3,  1, 0, 0, 2, 15),
  @(27, "Undrafted", "No",  "Charles Bediako",       "ALA", "MIZ@ALA", "9:52 - 1st Half",  10, 6,  4,  0, 1, 0, 0, 10),
  @(28, "Undrafted", "No",  "Shawn Phillips Jr.",    "MIZ", "MIZ@ALA", "9:52 - 1st Half",  9,  8,  5,  0, 0, 0, 1, 16),
  @(29, "Undrafted", "No",  "London Jemison",        "ALA", "MIZ@ALA", "9:52 - 1st Half",  7,  6,  2,  0, 0, 0, 0, 9),
  @(30, "Undrafted", "No",  "T.O. Barrett",          "MIZ", "MIZ@ALA", "9:52 - 1st Half",  4,  6,  1,  2, 0, 0, 2, 15),
  @(31, "Undrafted", "No",  "Trent Pierce",          "MIZ", "MIZ@ALA", "9:52 - 1st Half",  2,  2,  3,  0, 0, 0, 0, 9),
  @(32, "Undrafted", "No",  "D.J. Wagner",           "ARK", "ARK@OU",  "11:16 - 2nd Half", 0,  0,  2,  0, 0, 0, 1, 12),
  @(33, "Undrafted", "No",  "Nicholas Randall",      "MIZ", "MIZ@ALA", "9:52 - 1st Half",  0,  0,  0,  0, 0, 0, 0, 1),
  @(34, "Undrafted", "No",  "Sebastian Mack",        "MIZ", "MIZ@ALA", "9:52 - 1st Half",  -1, 0,  0,  0, 0, 0, 0, 0),
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  for ($col = 2; $col -le 15; $col++) {
    $ws1.Cells.Item($rowNum, $col).Value = $r[$col - 1]
  }
}

# Row 34 is brand new -- column A (date) needs to be written explicitly,
# forced to stay text (same "YYYY-MM-DD" string format used by every other
# row) instead of being auto-converted to a date serial number.
$ws1.Cells.Item(34, 1).NumberFormat = "@"
$ws1.Cells.Item(34, 1).Value = "2026-01-27"

# OwnerTotals sheet: recomputed starter_pooh_total / starters_count_so_far,
# re-sorted descending by starter_pooh_total (row 2, "Mark", is unchanged).
$ownerRows = @(
  @(2, "Mark", 21, 1),
  @(3, "Tar",  18, 2),
  @(4, "Booz", 13, 1),
  @(5, "CDL",  13, 1),
  @(6, "Clay", 9,  2),
  @(7, "Hal",  6,  1),
  @(8, "Ron",  0,  0),
)

foreach ($r in $ownerRows) {
  $rowNum = $r[0]
  $ws2.Cells.Item($rowNum, 1).Value = $r[1]
  $ws2.Cells.Item($rowNum, 2).Value = $r[2]
  $ws2.Cells.Item($rowNum, 3).Value = $r[3]
}
